# Add two new columns, I ("I0") and J ("IF"), matching the header style
# already used by the other header cells (A1:H1), then fill in the data
# rows 2-9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new headers "I0" and "IF" in I1/J1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Clone the existing header formatting (bold, bordered, centered) from H1
# onto the new header cells so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows: column I mirrors a constant "1" for rows 2-8, and column J
# mirrors the H column's value for those same rows.
$hValues = @{
    2 = 4
    3 = 5
    4 = 6
    5 = 6
    6 = 6
    7 = 5
    8 = 4
}

foreach ($row in $hValues.Keys) {
    $ws.Cells.Item($row, 9).Value = 1
    $ws.Cells.Item($row, 10).Value = $hValues[$row]
}

# Row 9 breaks the pattern seen in rows 2-8.
$ws.Cells.Item(9, 9).Value = 4
$ws.Cells.Item(9, 10).Value = 5
